$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Price" column cells are stored as plain text in the source sheet (values
# like "23.959.11" or "1.002" are dotted/rounded display strings, not real
# numbers). Writing a numeric-looking string straight into `.Value` lets Excel
# auto-convert it to a real number, so for those rows we prefix the literal
# with an apostrophe (forces text entry, just like a user would) and then
# restore the cell's original `.Style` afterwards so no formatting is left
# behind by the transient quote-prefix flag.

$ws.Range("D2").Value = "23.956.48"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "1.653.27"
$ws.Range("E3").Value = "  +2.15%  "

$s = $ws.Range("D4").Style
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = $s
$ws.Range("E4").Value = "  +0.00%  "

$s = $ws.Range("D5").Style
$ws.Range("D5").Value = "'309.67"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  +0.40%  "

$s = $ws.Range("D6").Style
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = $s
$ws.Range("E6").Value = "  -0.03%  "

$s = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.3906"
$ws.Range("D7").Style = $s
$ws.Range("E7").Value = "  -0.97%  "

$s = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.3838"
$ws.Range("D8").Style = $s
$ws.Range("E8").Value = "  +0.08%  "

$s = $ws.Range("D9").Style
$ws.Range("D9").Value = "'51.30"
$ws.Range("D9").Style = $s
$ws.Range("E9").Value = "  +3.96%  "

$s = $ws.Range("D10").Style
$ws.Range("D10").Value = "'1.356"
$ws.Range("D10").Style = $s
$ws.Range("E10").Value = "  +0.08%  "

$ws.Range("E11").Value = "  +0.01%  "

$s = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.08453"
$ws.Range("D12").Style = $s
$ws.Range("E12").Value = "  +0.15%  "

$s = $ws.Range("D13").Style
$ws.Range("D13").Value = "'23.99"
$ws.Range("D13").Style = $s
$ws.Range("E13").Value = "  +1.29%  "

$s = $ws.Range("D14").Style
$ws.Range("D14").Value = "'7.120"
$ws.Range("D14").Style = $s
$ws.Range("E14").Value = "  +1.27%  "

$s = $ws.Range("D15").Style
$ws.Range("D15").Value = "'7.887"
$ws.Range("D15").Style = $s
$ws.Range("E15").Value = "  +4.05%  "

$ws.Range("E16").Value = "  +3.08%  "

$ws.Range("D17").Value = "1.654.63"
$ws.Range("E17").Value = "  +2.40%  "

$s = $ws.Range("D18").Style
$ws.Range("D18").Value = "'94.66"
$ws.Range("D18").Style = $s

$s = $ws.Range("D19").Style
$ws.Range("D19").Value = "'0.06994"
$ws.Range("D19").Style = $s
$ws.Range("E19").Value = "  +0.93%  "

$s = $ws.Range("D20").Style
$ws.Range("D20").Value = "'19.77"
$ws.Range("D20").Style = $s
$ws.Range("E20").Value = "  -0.91%  "

$s = $ws.Range("D21").Style
$ws.Range("D21").Value = "'6.930"
$ws.Range("D21").Style = $s
$ws.Range("E21").Value = "  +1.81%  "

$s = $ws.Range("D22").Style
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = "  -0.07%  "

$s = $ws.Range("D23").Style
$ws.Range("D23").Value = "'13.66"
$ws.Range("D23").Style = $s
$ws.Range("E23").Value = "  +1.77%  "

$ws.Range("D24").Value = "23.964.29"
$ws.Range("E24").Value = "  +0.53%  "

$s = $ws.Range("D25").Style
$ws.Range("D25").Value = "'2.489"
$ws.Range("D25").Style = $s
$ws.Range("E25").Value = "  +1.97%  "

$s = $ws.Range("D26").Style
$ws.Range("D26").Value = "'3.013"
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = "  +6.36%  "

$ws.Range("E27").Value = "  -0.33%  "

$s = $ws.Range("D28").Style
$ws.Range("D28").Value = "'151.19"
$ws.Range("D28").Style = $s
$ws.Range("E28").Value = "  -3.68%  "

$s = $ws.Range("D29").Style
$ws.Range("D29").Value = "'5.448"
$ws.Range("D29").Style = $s
$ws.Range("E29").Value = "  +3.17%  "

$s = $ws.Range("D30").Style
$ws.Range("D30").Value = "'139.52"
$ws.Range("D30").Style = $s
$ws.Range("E30").Value = "  -0.25%  "

$s = $ws.Range("D31").Style
$ws.Range("D31").Value = "'7.888"
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = "  +1.04%  "

$s = $ws.Range("D32").Style
$ws.Range("D32").Value = "'2.490"
$ws.Range("D32").Style = $s
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").Value = "1.835.69"
$ws.Range("E33").Value = "  +2.30%  "

$ws.Range("E34").Value = "  +7.04%  "

$s = $ws.Range("D35").Style
$ws.Range("D35").Value = "'0.08116"
$ws.Range("D35").Style = $s
$ws.Range("E35").Value = "  +0.29%  "

$s = $ws.Range("D36").Style
$ws.Range("D36").Value = "'0.02971"
$ws.Range("D36").Style = $s
$ws.Range("E36").Value = "  +3.33%  "

$s = $ws.Range("D37").Style
$ws.Range("D37").Value = "'6.778"
$ws.Range("D37").Style = $s
$ws.Range("E37").Value = "  +3.08%  "

$s = $ws.Range("D38").Style
$ws.Range("D38").Value = "'10.87"
$ws.Range("D38").Style = $s
$ws.Range("E38").Value = "  +5.18%  "

$s = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.2684"
$ws.Range("D39").Style = $s
$ws.Range("E39").Value = "  +0.83%  "

$s = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.09140"
$ws.Range("D40").Style = $s
$ws.Range("E40").Value = "  +0.09%  "

$s = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.7563"
$ws.Range("D41").Style = $s
$ws.Range("E41").Value = "  +0.96%  "

$s = $ws.Range("D42").Style
$ws.Range("D42").Value = "'13.45"
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = "  -0.54%  "

$s = $ws.Range("D43").Style
$ws.Range("D43").Value = "'1.425"
$ws.Range("D43").Style = $s
$ws.Range("E43").Value = "  +0.02%  "

$s = $ws.Range("D44").Style
$ws.Range("D44").Value = "'16.41"
$ws.Range("D44").Style = $s
$ws.Range("E44").Value = "  +1.56%  "

$s = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.6948"
$ws.Range("D45").Style = $s
$ws.Range("E45").Value = "  +0.56%  "

$s = $ws.Range("D46").Style
$ws.Range("D46").Value = "'2.457"
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = "  -0.55%  "

$s = $ws.Range("D47").Style
$ws.Range("D47").Value = "'4.089"
$ws.Range("D47").Style = $s
$ws.Range("E47").Value = "  +0.47%  "

$s = $ws.Range("D48").Style
$ws.Range("D48").Value = "'1.000"
$ws.Range("D48").Style = $s
$ws.Range("E48").Value = "  -0.06%  "

$s = $ws.Range("D49").Style
$ws.Range("D49").Value = "'0.08287"
$ws.Range("D49").Style = $s
$ws.Range("E49").Value = "  +0.84%  "

$s = $ws.Range("D50").Style
$ws.Range("D50").Value = "'134.60"
$ws.Range("D50").Style = $s
$ws.Range("E50").Value = "  +0.09%  "

$s = $ws.Range("D51").Style
$ws.Range("D51").Value = "'1.208"
$ws.Range("D51").Style = $s
$ws.Range("E51").Value = "  +0.49%  "

